$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (stored OOXML width = ColumnWidth + 5/6) ---
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 6.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 6.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(14).ColumnWidth = 6.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(18).ColumnWidth = 6.166666666666667
$ws.Columns.Item(19).ColumnWidth = 6.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(25).ColumnWidth = 6.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(27).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(31).ColumnWidth = 6.166666666666667
$ws.Columns.Item(32).ColumnWidth = 7.166666666666667
$ws.Columns.Item(33).ColumnWidth = 6.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667

# --- Update data rows 2-5 with new values ---
$ws.Cells.Item(2, 1).Value = 45110.50694444445
$ws.Cells.Item(2, 2).Value = 10.232
$ws.Cells.Item(2, 3).Value = 6.988
$ws.Cells.Item(2, 4).Value = 3.382
$ws.Cells.Item(2, 5).Value = 22.493
$ws.Cells.Item(2, 6).Value = 16.551
$ws.Cells.Item(2, 7).Value = 7.711
$ws.Cells.Item(2, 8).Value = 22.993
$ws.Cells.Item(2, 9).Value = 12.747
$ws.Cells.Item(2, 10).Value = 5.033
$ws.Cells.Item(2, 11).Value = 6.872
$ws.Cells.Item(2, 12).Value = 8.682
$ws.Cells.Item(2, 13).Value = 9.968
$ws.Cells.Item(2, 14).Value = 2.214
$ws.Cells.Item(2, 15).Value = 8.271
$ws.Cells.Item(2, 16).Value = 11.091
$ws.Cells.Item(2, 17).Value = 7.662
$ws.Cells.Item(2, 18).Value = 2.404
$ws.Cells.Item(2, 19).Value = 1.071
$ws.Cells.Item(2, 20).Value = 118.497
$ws.Cells.Item(2, 21).Value = 22.775
$ws.Cells.Item(2, 22).Value = 7.635
$ws.Cells.Item(2, 23).Value = 14.429
$ws.Cells.Item(2, 24).Value = 7.663
$ws.Cells.Item(2, 25).Value = 2.148
$ws.Cells.Item(2, 26).Value = 12.965
$ws.Cells.Item(2, 27).Value = 6.744
$ws.Cells.Item(2, 28).Value = 6.372
$ws.Cells.Item(2, 29).Value = 7.242
$ws.Cells.Item(2, 30).Value = 9.571
$ws.Cells.Item(2, 31).Value = 2.474
$ws.Cells.Item(2, 32).Value = 20.428
$ws.Cells.Item(2, 33).Value = 3.858
$ws.Cells.Item(2, 34).Value = 9.544

$ws.Cells.Item(3, 1).Value = 45110.51388888889
$ws.Cells.Item(3, 2).Value = 17.552
$ws.Cells.Item(3, 3).Value = 12.927
$ws.Cells.Item(3, 4).Value = 1.796
$ws.Cells.Item(3, 5).Value = 38.574
$ws.Cells.Item(3, 6).Value = 30.767
$ws.Cells.Item(3, 7).Value = 13.626
$ws.Cells.Item(3, 8).Value = 51.053
$ws.Cells.Item(3, 9).Value = 21.477
$ws.Cells.Item(3, 10).Value = 9.445
$ws.Cells.Item(3, 11).Value = 13.487
$ws.Cells.Item(3, 12).Value = 15.368
$ws.Cells.Item(3, 13).Value = 16.636
$ws.Cells.Item(3, 14).Value = 4.204
$ws.Cells.Item(3, 15).Value = 13.911
$ws.Cells.Item(3, 16).Value = 19.536
$ws.Cells.Item(3, 17).Value = 12.065
$ws.Cells.Item(3, 18).Value = 1.26
$ws.Cells.Item(3, 19).Value = 0.799
$ws.Cells.Item(3, 20).Value = 204.407
$ws.Cells.Item(3, 21).Value = 38.856
$ws.Cells.Item(3, 22).Value = 12.84
$ws.Cells.Item(3, 23).Value = 25.771
$ws.Cells.Item(3, 24).Value = 13.541
$ws.Cells.Item(3, 25).Value = 2.234
$ws.Cells.Item(3, 26).Value = 25.96
$ws.Cells.Item(3, 27).Value = 11.342
$ws.Cells.Item(3, 28).Value = 10.227
$ws.Cells.Item(3, 29).Value = 11.962
$ws.Cells.Item(3, 30).Value = 16.278
$ws.Cells.Item(3, 31).Value = 1.101
$ws.Cells.Item(3, 32).Value = 46.553
$ws.Cells.Item(3, 33).Value = 7.039
$ws.Cells.Item(3, 34).Value = 16.053

$ws.Cells.Item(4, 1).Value = 45110.52083333334
$ws.Cells.Item(4, 2).Value = 13.765
$ws.Cells.Item(4, 3).Value = 10.199
$ws.Cells.Item(4, 4).Value = 1.215
$ws.Cells.Item(4, 5).Value = 30.269
$ws.Cells.Item(4, 6).Value = 24.265
$ws.Cells.Item(4, 7).Value = 10.699
$ws.Cells.Item(4, 8).Value = 43.741
$ws.Cells.Item(4, 9).Value = 16.826
$ws.Cells.Item(4, 10).Value = 7.486
$ws.Cells.Item(4, 11).Value = 10.652
$ws.Cells.Item(4, 12).Value = 12.093
$ws.Cells.Item(4, 13).Value = 13.038
$ws.Cells.Item(4, 14).Value = 3.315
$ws.Cells.Item(4, 15).Value = 10.903
$ws.Cells.Item(4, 16).Value = 15.357
$ws.Cells.Item(4, 17).Value = 9.42
$ws.Cells.Item(4, 18).Value = 0.869
$ws.Cells.Item(4, 19).Value = 0.561
$ws.Cells.Item(4, 20).Value = 158.636
$ws.Cells.Item(4, 21).Value = 30.536
$ws.Cells.Item(4, 22).Value = 10.064
$ws.Cells.Item(4, 23).Value = 20.298
$ws.Cells.Item(4, 24).Value = 10.646
$ws.Cells.Item(4, 25).Value = 1.683
$ws.Cells.Item(4, 26).Value = 21.447
$ws.Cells.Item(4, 27).Value = 8.889
$ws.Cells.Item(4, 28).Value = 7.994
$ws.Cells.Item(4, 29).Value = 9.365
$ws.Cells.Item(4, 30).Value = 12.77
$ws.Cells.Item(4, 31).Value = 0.708
$ws.Cells.Item(4, 32).Value = 39.873
$ws.Cells.Item(4, 33).Value = 5.541
$ws.Cells.Item(4, 34).Value = 12.582

$ws.Cells.Item(5, 1).Value = 45110.52777777778
$ws.Cells.Item(5, 2).Value = 23.41
$ws.Cells.Item(5, 3).Value = 17.49
$ws.Cells.Item(5, 4).Value = 1.33
$ws.Cells.Item(5, 5).Value = 51.19
$ws.Cells.Item(5, 6).Value = 41.75
$ws.Cells.Item(5, 7).Value = 18.31
$ws.Cells.Item(5, 8).Value = 70.51
$ws.Cells.Item(5, 9).Value = 28.46
$ws.Cells.Item(5, 10).Value = 12.76
$ws.Cells.Item(5, 11).Value = 18.61
$ws.Cells.Item(5, 12).Value = 20.51
$ws.Cells.Item(5, 13).Value = 21.88
$ws.Cells.Item(5, 14).Value = 5.77
$ws.Cells.Item(5, 15).Value = 18.42
$ws.Cells.Item(5, 16).Value = 26.16
$ws.Cells.Item(5, 17).Value = 15.56
$ws.Cells.Item(5, 18).Value = 0.77
$ws.Cells.Item(5, 19).Value = 0.79
$ws.Cells.Item(5, 20).Value = 273.09
$ws.Cells.Item(5, 21).Value = 51.49
$ws.Cells.Item(5, 22).Value = 17
$ws.Cells.Item(5, 23).Value = 34.59
$ws.Cells.Item(5, 24).Value = 18.16
$ws.Cells.Item(5, 25).Value = 2.6
$ws.Cells.Item(5, 26).Value = 34.79
$ws.Cells.Item(5, 27).Value = 15.02
$ws.Cells.Item(5, 28).Value = 13.33
$ws.Cells.Item(5, 29).Value = 15.65
$ws.Cells.Item(5, 30).Value = 21.6
$ws.Cells.Item(5, 31).Value = 0.52
$ws.Cells.Item(5, 32).Value = 64.01
$ws.Cells.Item(5, 33).Value = 9.55
$ws.Cells.Item(5, 34).Value = 21.26

# --- Delete row 6 (data reduced from 5 rows to 4 rows) ---
$ws.Rows.Item(6).Delete()